# Update the LR-pair table with the new TPM-derived values.
# The sending/target cluster grid is now a full 3x3 (ECs/FAPs/MuSCs x ECs/FAPs/MuSCs),
# replacing the old 6-row (2x3) table that omitted the ECs target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",   "Ccl2", "Ackr4", "ECs",   2, 0.6666666666666666, 0.1510096666666667, 0.453029,           0.01105950042918124, 0.01105950042918124, 1, 0.3333333333333333, 0.07058066666666667, 0.211742, 0.03405243646146196, 0.03405243646146196, 0.01065836294644444, 0.09592526651800001, 0.0003766029356602054, 0.0003766029356602055),
    @("ECs",   "Ccl2", "Ackr4", "FAPs",  2, 0.6666666666666666, 0.1510096666666667, 0.453029,           0.01105950042918124, 0.01105950042918124, 3, 1,                   1.890486333333333,   5.671459, 0.9120863940138783,  0.9120863940138783,  0.2854817110345556,  2.569335399311,      0.01008721986604686,   0.01008721986604686),
    @("ECs",   "Ccl2", "Ackr4", "MuSCs", 2, 0.6666666666666666, 0.1510096666666667, 0.453029,           0.01105950042918124, 0.01105950042918124, 2, 0.6666666666666666, 0.1116383333333333,  0.334915, 0.05386116952465989, 0.05386116952465989, 0.01685846750388889, 0.151726207535,      0.0005956776274741796, 0.0005956776274741796),
    @("FAPs",  "Ccl2", "Ackr4", "ECs",   3, 1,                   10.23495333333333, 30.70486,           0.7495776481151314,  0.7495776481151314,  1, 0.3333333333333333, 0.07058066666666667, 0.211742, 0.03405243646146196, 0.03405243646146196, 0.7223898295688889,  6.501508466120001,  0.0255249452353726,    0.02552494523537261),
    @("FAPs",  "Ccl2", "Ackr4", "FAPs",  3, 1,                   10.23495333333333, 30.70486,           0.7495776481151314,  0.7495776481151314,  3, 1,                   1.890486333333333,   5.671459, 0.9120863940138783,  0.9120863940138783,  19.34903939897111,   174.14135459074,    0.6836795741027339,    0.6836795741027339),
    @("FAPs",  "Ccl2", "Ackr4", "MuSCs", 3, 1,                   10.23495333333333, 30.70486,           0.7495776481151314,  0.7495776481151314,  2, 0.6666666666666666, 0.1116383333333333,  0.334915, 0.05386116952465989, 0.05386116952465989, 1.142613131877778,   10.2835181869,      0.04037312877702495,   0.04037312877702495),
    @("MuSCs", "Ccl2", "Ackr4", "ECs",   3, 1,                   3.268330666666666, 9.804991999999999, 0.2393628514556874,  0.2393628514556874,  1, 0.3333333333333333, 0.07058066666666667, 0.211742, 0.03405243646146196, 0.03405243646146196, 0.2306809573404444,  2.076128616064,     0.008150888290429151,  0.008150888290429153),
    @("MuSCs", "Ccl2", "Ackr4", "FAPs",  3, 1,                   3.268330666666666, 9.804991999999999, 0.2393628514556874,  0.2393628514556874,  3, 1,                   1.890486333333333,   5.671459, 0.9120863940138783,  0.9120863940138783,  6.178734458147555,   55.608610123328,    0.2183196000450975,    0.2183196000450975),
    @("MuSCs", "Ccl2", "Ackr4", "MuSCs", 3, 1,                   3.268330666666666, 9.804991999999999, 0.2393628514556874,  0.2393628514556874,  2, 0.6666666666666666, 0.1116383333333333,  0.334915, 0.05386116952465989, 0.05386116952465989, 0.3648709884088888,  3.28383889568,      0.01289236312016076,   0.01289236312016076)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
